# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-27 12:29:12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I (9th column) width: 14 -> 10 ---
# Excel's ColumnWidth property is offset from the stored OOXML <col width> by
# ~0.8333 (5/6) characters of padding, so 9.1666.. yields a stored width of 10.
$ws.Columns.Item(9).ColumnWidth = 9.1666666666667

# --- Top summary block (K/L columns) ---
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 0

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "40.5%"

# --- Per-group breakdown table (rows 15-26), columns O/P/R/S ---
function Set-PctText($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

$ws.Range("O15").Value = 7
$ws.Range("P15").Value = 0
Set-PctText "R15" "38.9%"
Set-PctText "S15" "76.7%"

$ws.Range("O16").Value = 8
$ws.Range("P16").Value = 0
Set-PctText "R16" "42.1%"
Set-PctText "S16" "71.8%"

$ws.Range("O17").Value = 8
$ws.Range("P17").Value = 0
Set-PctText "R17" "42.1%"
Set-PctText "S17" "59.9%"

$ws.Range("O18").Value = 8
$ws.Range("P18").Value = 0
Set-PctText "R18" "42.1%"
Set-PctText "S18" "83.9%"

$ws.Range("O19").Value = 7
$ws.Range("P19").Value = 0
Set-PctText "R19" "38.9%"
Set-PctText "S19" "88.0%"

$ws.Range("O20").Value = 7
$ws.Range("P20").Value = 0
Set-PctText "R20" "38.9%"
Set-PctText "S20" "90.8%"

$ws.Range("O21").Value = 7
$ws.Range("P21").Value = 0
Set-PctText "R21" "38.9%"
Set-PctText "S21" "89.7%"

$ws.Range("O22").Value = 7
$ws.Range("P22").Value = 0
Set-PctText "R22" "38.9%"
Set-PctText "S22" "90.5%"

$ws.Range("O23").Value = 7
$ws.Range("P23").Value = 0
Set-PctText "R23" "38.9%"
Set-PctText "S23" "67.6%"

$ws.Range("O24").Value = 8
$ws.Range("P24").Value = 0
Set-PctText "R24" "42.1%"
Set-PctText "S24" "70.4%"

$ws.Range("O25").Value = 8
$ws.Range("P25").Value = 0
Set-PctText "R25" "42.1%"
Set-PctText "S25" "74.6%"

$ws.Range("O26").Value = 8
$ws.Range("P26").Value = 0
Set-PctText "R26" "42.1%"
Set-PctText "S26" "71.1%"

# --- Session rows that flip from "Not Recorded" (pink) to "Recorded" (green) ---
# Re-style A:I on each row to match the "Recorded" look (copy format from a
# known-good Recorded row, A2:I2), then set the Recorded-By/Students/Status cells.
$styleSource = $ws.Range("A2:I2")

function Set-RecordedRow($row, $studentsText) {
    $dst = $ws.Range("A" + $row + ":I" + $row)
    $styleSource.Copy()
    $dst.PasteSpecial(-4122)
    $ws.Range("G" + $row).Value = "dnasr281@gmail.com"
    $ws.Range("H" + $row).Value = $studentsText
    $ws.Range("I" + $row).Value = "Recorded"
}

Set-RecordedRow 16  "21/27"
Set-RecordedRow 35  "18/31"
Set-RecordedRow 54  "13/19"
Set-RecordedRow 73  "18/21"
Set-RecordedRow 91  "25/31"
Set-RecordedRow 109 "25/28"
Set-RecordedRow 127 "21/29"
Set-RecordedRow 145 "28/33"
Set-RecordedRow 163 "24/30"
Set-RecordedRow 182 "20/27"
Set-RecordedRow 201 "23/29"
Set-RecordedRow 220 "24/29"

$excel.CutCopyMode = $false
